# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which carry identical data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2 = 360
    3 = 99
    4 = 1551
    6 = 54
    9 = 395
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
